# Cost-benefit worksheet update: bump the Codenomicon "Implementation - Licensing
# (Year 1)" figure (F3) from 10000 to 15000, let dependent formulas recalc, and
# leave the selection on F4 as the final active cell (matching the saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 15000

$excel.Calculate()

$ws.Range("F4").Select()
